$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin names & links) ---
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# --- Numeric-looking cells stored as text (Price / Volume columns) ---
# Mark these ranges as Text format first so Excel keeps the literal
# string (matching the original inlineStr cells) instead of coercing
# to a Number/Percentage value.
$numRanges = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","E18","E19","E20","D21","E21","E22","D23","E23","D24","E24","D25","E25","D26","E26","E27","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $numRanges) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "303.69"
$ws.Range("E2").Value = "-0.10%"
$ws.Range("D3").Value = "37.34"
$ws.Range("E3").Value = "4.17%"
$ws.Range("D4").Value = "5.032"
$ws.Range("E4").Value = "-0.96%"
$ws.Range("D5").Value = "0.07839"
$ws.Range("E5").Value = "-0.38%"
$ws.Range("D6").Value = "2.208"
$ws.Range("E6").Value = "-3.77%"
$ws.Range("D7").Value = "7.985"
$ws.Range("E7").Value = "-1.16%"
$ws.Range("D8").Value = "4.027"
$ws.Range("E8").Value = "0.79%"
$ws.Range("D9").Value = "0.9256"
$ws.Range("E9").Value = "0.04%"
$ws.Range("D10").Value = "0.09875"
$ws.Range("E10").Value = "-1.66%"
$ws.Range("D11").Value = "0.1882"
$ws.Range("E11").Value = "3.04%"
$ws.Range("D12").Value = "0.08668"
$ws.Range("E12").Value = "0.26%"
$ws.Range("D13").Value = "0.03617"
$ws.Range("E13").Value = "6.61%"
$ws.Range("D14").Value = "0.09942"
$ws.Range("E14").Value = "0.25%"
$ws.Range("D15").Value = "0.001479"
$ws.Range("E15").Value = "-0.19%"
$ws.Range("D16").Value = "0.005690"
$ws.Range("E16").Value = "1.61%"
$ws.Range("D17").Value = "3.457"
$ws.Range("E17").Value = "-0.89%"
$ws.Range("E18").Value = "12.00%"
$ws.Range("E19").Value = "0.01%"
$ws.Range("E20").Value = "0.61%"
$ws.Range("D21").Value = "4.774"
$ws.Range("E21").Value = "4.73%"
$ws.Range("E22").Value = "-0.81%"
$ws.Range("D23").Value = "0.04606"
$ws.Range("E23").Value = "-1.19%"
$ws.Range("D24").Value = "0.005198"
$ws.Range("E24").Value = "15.77%"
$ws.Range("D25").Value = "0.001252"
$ws.Range("E25").Value = "0.97%"
$ws.Range("D26").Value = "0.0001401"
$ws.Range("E26").Value = "7.86%"
$ws.Range("E27").Value = "-9.27%"
$ws.Range("D39").Value = "0.01829"
$ws.Range("E39").Value = "4.55%"
$ws.Range("D40").Value = "0.04751"
$ws.Range("E40").Value = "1.21%"
$ws.Range("D41").Value = "0.007940"
$ws.Range("E41").Value = "1.27%"
$ws.Range("D42").Value = "0.1402"
$ws.Range("E42").Value = "-1.13%"
$ws.Range("D43").Value = "0.007567"
$ws.Range("E43").Value = "-13.98%"
$ws.Range("D44").Value = "0.002242"
$ws.Range("E44").Value = "1.49%"
$ws.Range("D45").Value = "0.01039"
$ws.Range("E45").Value = "13.13%"
$ws.Range("D46").Value = "0.00006304"
$ws.Range("E46").Value = "5.08%"
$ws.Range("E47").Value = "0.15%"
$ws.Range("D48").Value = "0.0005804"
$ws.Range("E48").Value = "0.06%"
$ws.Range("D49").Value = "35.90"
$ws.Range("E49").Value = "518.98%"
$ws.Range("D50").Value = "0.002691"
$ws.Range("E50").Value = "0.14%"
$ws.Range("D51").Value = "0.00002102"
$ws.Range("E51").Value = "0.15%"
